# report912 / sum912.xlsx : update the energy-equivalent (ktoe) columns
# (AM/AO/AQ) in table 9.7 so that they are derived using the new
# kWh/kg -> ktoe conversion factors that now come from the settings
# database instead of being baked into the worksheet.
#
#   - electric appliances (rows 13-31): old factor 0.08521   -> new factor 0.09
#   - gas stoves           (rows 32-34): old factor 0.0000103194 -> new factor 0.0000084
#   - solid biomass fuels   (rows 35-52): old factor 0.000000004334148 -> new factor 0.000000003528
#
# Each "quantity" column (AL, AN, AP) keeps its physical-unit value; the
# paired "ktoe" column (AM, AO, AQ respectively) is simply that quantity
# multiplied by the conversion factor for its fuel/device group, so we
# recompute AM/AO/AQ from the existing AL/AN/AP values using the new
# factors rather than hard-coding each result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newElectricFactor = 0.09
$newGasFactor = 0.0000084
$newBiomassFactor = 0.000000003528

$quantityCols = "AL", "AN", "AP"
$ktoeCols = "AM", "AO", "AQ"

for ($row = 13; $row -le 52; $row++) {
    if ($row -ge 13 -and $row -le 31) {
        $factor = $newElectricFactor
    } elseif ($row -ge 32 -and $row -le 34) {
        $factor = $newGasFactor
    } else {
        $factor = $newBiomassFactor
    }

    for ($i = 0; $i -lt 3; $i++) {
        $qCell = $quantityCols[$i] + $row
        $kCell = $ktoeCols[$i] + $row
        $qty = $ws.Range($qCell).Value2
        $ws.Range($kCell).Value2 = $qty * $factor
    }
}
